$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (id) and C (phone) hold digit-only strings that Excel would
# otherwise auto-coerce into numbers (scientific notation, lost leading
# zeroes, etc). Format those cells as Text before writing so the values
# are kept verbatim as strings, then clear the formatting again so no
# extra cell style is left behind (source file has none).
$ws.Range("A1:A3").NumberFormat = "@"
$ws.Range("C1:C3").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "designation"
$ws.Range("F1").Value = "photo"
$ws.Range("G1").Value = "createdAt"

# Row 2
$ws.Range("A2").Value = "1752670313518"
$ws.Range("B2").Value = "Abu Inshah"
$ws.Range("C2").Value = "7449085120"
$ws.Range("D2").Value = "aiautomationhig@gmail.com"
$ws.Range("E2").Value = "Health insurance advisor,Wealth Manager"
$ws.Range("F2").Value = "uploads/abu_inshah_1752670313452.jpeg"
$ws.Range("G2").Value = "2025-07-16T12:51:53.518Z"

# Row 3
$ws.Range("A3").Value = "1752671851601"
$ws.Range("B3").Value = "Abu Inshah"
$ws.Range("C3").Value = "7449085120"
$ws.Range("D3").Value = "ajai17101999@gmail.com"
$ws.Range("E3").Value = "Health insurance advisor,Wealth Manager"
$ws.Range("F3").Value = "uploads/abu_inshah_1752671851576.jpeg"
$ws.Range("G3").Value = "2025-07-16T13:17:31.601Z"

# The Text number format was only needed to stop Excel re-typing the
# digit-only strings as numbers while the values were entered; strip it
# again now that the text values are locked in, so no stray cell style
# is left on the sheet.
$ws.Range("A1:A3").ClearFormats()
$ws.Range("C1:C3").ClearFormats()
